$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet so the cells can be updated
$ws.Unprotect()

# Update the "as of" date in the confidentiality disclaimer
$dateCell = $ws.Cells.Find("2021-07-08")
if ($dateCell) {
    $oldText = $dateCell.Text
    $newText = $oldText -replace "2021-07-08", "2021-07-09"
    $dateCell.Value = $newText
} else {
    $oldText = $ws.Range("A40").Text
    $newText = $oldText -replace "2021-07-08", "2021-07-09"
    $ws.Range("A40").Value = $newText
}

# Update the Weight (D) and Percent Change (E) values for each holding row
$ws.Range("D2").Value = 0.03075688484814154
$ws.Range("E2").Value = 0.001976794155565065
$ws.Range("D3").Value = 0.03432463817373325
$ws.Range("E3").Value = -0.002392936710949867
$ws.Range("D4").Value = 0.03388940814380529
$ws.Range("E4").Value = 0.01292623133496762
$ws.Range("D5").Value = 0.07045647227650796
$ws.Range("E5").Value = -0.003234702163525283
$ws.Range("D6").Value = 0.03065681026425355
$ws.Range("E6").Value = -0.002069475240206842
$ws.Range("D7").Value = 0.01566714815758915
$ws.Range("E7").Value = 0.02270590787475602
$ws.Range("D8").Value = 0.03228745952171902
$ws.Range("E8").Value = 0.003766169968888011
$ws.Range("D9").Value = 0.03154086536192068
$ws.Range("E9").Value = 0.01496030938326887
$ws.Range("D10").Value = 0.050109043532444
$ws.Range("E10").Value = 0.009834953651367773
$ws.Range("D11").Value = 0.02641496964718989
$ws.Range("E11").Value = 0.03047285464098071
$ws.Range("D12").Value = 0.01514449448173639
$ws.Range("E12").Value = 0.01754232850410209
$ws.Range("D13").Value = 0.01599682782451072
$ws.Range("E13").Value = 0.02089235127478761
$ws.Range("D14").Value = 0.01508105097195079
$ws.Range("E14").Value = 0.009728308501314675
$ws.Range("D15").Value = 0.006529016908828279
$ws.Range("E15").Value = 0.06194690265486735
$ws.Range("D16").Value = 0.007055824624012236
$ws.Range("E16").Value = 0.04324555769642457
$ws.Range("D17").Value = 0.0324883639693734
$ws.Range("E17").Value = -0.003010577705451656
$ws.Range("D18").Value = 0.0275267416281946
$ws.Range("E18").Value = 0.01788274352976704
$ws.Range("D19").Value = 0.03111904154983431
$ws.Range("E19").Value = 0.00152905198776776
$ws.Range("D20").Value = 0.03263281124611739
$ws.Range("E20").Value = 0.01380008679299882
$ws.Range("D21").Value = 0.0487823944260345
$ws.Range("E21").Value = 0.003077173180984216
$ws.Range("D22").Value = 0.02711001595528743
$ws.Range("E22").Value = 0.03571627570067415
$ws.Range("D23").Value = 0.03007335655818959
$ws.Range("E23").Value = 0.01114459722483807
$ws.Range("D24").Value = 0.02758263233919619
$ws.Range("E24").Value = 0.02137884295483961
$ws.Range("D25").Value = 0.01200139726777504
$ws.Range("E25").Value = 0.03624921334172426
$ws.Range("D26").Value = 0.01265018268332059
$ws.Range("E26").Value = 0.01979222640157619
$ws.Range("D27").Value = 0.02923197477365206
$ws.Range("E27").Value = 0.00554213443228635
$ws.Range("D28").Value = 0.02848462533397532
$ws.Range("E28").Value = 0.00281062735323756
$ws.Range("D29").Value = 0.03142946157985668
$ws.Range("E29").Value = 0.001874414245548239
$ws.Range("D30").Value = 0.03336316688853013
$ws.Range("E30").Value = 0.002303430243416615
$ws.Range("D31").Value = 0.0300642931996488
$ws.Range("E31").Value = 0.007411036163344153
$ws.Range("D32").Value = 0.02829656064425374
$ws.Range("E32").Value = -0.002535699986654372
$ws.Range("D33").Value = 0.02893307276177529
$ws.Range("E33").Value = 0.02282175277848486
$ws.Range("D34").Value = 0.03096269861500552
$ws.Range("E34").Value = 0.005268935236004468
$ws.Range("D35").Value = 0.02850048621142172
$ws.Range("E35").Value = 0.03199946998807479
$ws.Range("D36").Value = 0.03285580763021498
$ws.Range("E36").Value = 0
$ws.Range("E37").Value = 0.00971384333607106

# Re-apply sheet protection to restore the original protected state
$ws.Protect()
